$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add manually-curated OC and BC SNPs for prostate cancer PRS (rows 31-72)
# Data is written in the same column-by-column order the original author used,
# so that the shared-string table is populated in a matching sequence.

# rows 31-45: rs_id column (ovarian cancer risk SNPs)
$ws.Range("A31").Value = "rs58722170"
$ws.Range("A32").Value = "rs711830"
$ws.Range("A33").Value = "rs62274041"
$ws.Range("A34").Value = "rs10069690"
$ws.Range("A35").Value = "rs78724141"
$ws.Range("A36").Value = "rs10088218"
$ws.Range("A37").Value = "rs7032221"
$ws.Range("A38").Value = "rs635634"
$ws.Range("A39").Value = "rs1802669"
$ws.Range("A40").Value = "rs7135337"
$ws.Range("A41").Value = "rs11651755"
$ws.Range("A42").Value = "rs1105569"
$ws.Range("A43").Value = "rs7207826"
$ws.Range("A44").Value = "rs61494113"
$ws.Range("A45").Value = "rs9625477"

# rows 31-45: notes column
$ws.Range("C31").Value = "ovarian cancer (Comprehensive epithelial tubo-ovarian cancer risk prediction model incorporating genetic and epidemiological risk factors; PMID: 34844974; Table S5)"
$ws.Range("C32").Value = "ovarian cancer (Comprehensive epithelial tubo-ovarian cancer risk prediction model incorporating genetic and epidemiological risk factors; PMID: 34844974; Table S5)"
$ws.Range("C33").Value = "ovarian cancer (Comprehensive epithelial tubo-ovarian cancer risk prediction model incorporating genetic and epidemiological risk factors; PMID: 34844974; Table S5)"
$ws.Range("C34").Value = "ovarian cancer (Comprehensive epithelial tubo-ovarian cancer risk prediction model incorporating genetic and epidemiological risk factors; PMID: 34844974; Table S5)"
$ws.Range("C35").Value = "ovarian cancer (Comprehensive epithelial tubo-ovarian cancer risk prediction model incorporating genetic and epidemiological risk factors; PMID: 34844974; Table S5)"
$ws.Range("C36").Value = "ovarian cancer (Comprehensive epithelial tubo-ovarian cancer risk prediction model incorporating genetic and epidemiological risk factors; PMID: 34844974; Table S5)"
$ws.Range("C37").Value = "ovarian cancer (Comprehensive epithelial tubo-ovarian cancer risk prediction model incorporating genetic and epidemiological risk factors; PMID: 34844974; Table S5)"
$ws.Range("C38").Value = "ovarian cancer (Comprehensive epithelial tubo-ovarian cancer risk prediction model incorporating genetic and epidemiological risk factors; PMID: 34844974; Table S5)"
$ws.Range("C39").Value = "ovarian cancer (Comprehensive epithelial tubo-ovarian cancer risk prediction model incorporating genetic and epidemiological risk factors; PMID: 34844974; Table S5)"
$ws.Range("C40").Value = "ovarian cancer (Comprehensive epithelial tubo-ovarian cancer risk prediction model incorporating genetic and epidemiological risk factors; PMID: 34844974; Table S5)"
$ws.Range("C41").Value = "ovarian cancer (Comprehensive epithelial tubo-ovarian cancer risk prediction model incorporating genetic and epidemiological risk factors; PMID: 34844974; Table S5)"
$ws.Range("C42").Value = "ovarian cancer (Comprehensive epithelial tubo-ovarian cancer risk prediction model incorporating genetic and epidemiological risk factors; PMID: 34844974; Table S5)"
$ws.Range("C43").Value = "ovarian cancer (Comprehensive epithelial tubo-ovarian cancer risk prediction model incorporating genetic and epidemiological risk factors; PMID: 34844974; Table S5)"
$ws.Range("C44").Value = "ovarian cancer (Comprehensive epithelial tubo-ovarian cancer risk prediction model incorporating genetic and epidemiological risk factors; PMID: 34844974; Table S5)"
$ws.Range("C45").Value = "ovarian cancer (Comprehensive epithelial tubo-ovarian cancer risk prediction model incorporating genetic and epidemiological risk factors; PMID: 34844974; Table S5)"

# rows 46-55: rs_id column (breast cancer risk SNPs - Table 5)
$ws.Range("A46").Value = "rs60882887"
$ws.Range("A47").Value = "rs5820435"
$ws.Range("A48").Value = "rs7222250"
$ws.Range("A49").Value = "rs9901834"
$ws.Range("A50").Value = "rs58117746"
$ws.Range("A51").Value = "rs2239711"
$ws.Range("A52").Value = "rs10708222"
$ws.Range("A53").Value = "rs41283425"
$ws.Range("A54").Value = "rs56291217"
$ws.Range("A55").Value = "rs111637825"

# rows 46-55: notes column (Table 5)
$ws.Range("C46").Value = "breast cancer (A case-only study to identify genetic modifiers of breast cancer risk for BRCA1/BRCA2 mutation carriers; Table 5; PMID: 33597508)"
$ws.Range("C47").Value = "breast cancer (A case-only study to identify genetic modifiers of breast cancer risk for BRCA1/BRCA2 mutation carriers; Table 5; PMID: 33597508)"
$ws.Range("C48").Value = "breast cancer (A case-only study to identify genetic modifiers of breast cancer risk for BRCA1/BRCA2 mutation carriers; Table 5; PMID: 33597508)"
$ws.Range("C49").Value = "breast cancer (A case-only study to identify genetic modifiers of breast cancer risk for BRCA1/BRCA2 mutation carriers; Table 5; PMID: 33597508)"
$ws.Range("C50").Value = "breast cancer (A case-only study to identify genetic modifiers of breast cancer risk for BRCA1/BRCA2 mutation carriers; Table 5; PMID: 33597508)"
$ws.Range("C51").Value = "breast cancer (A case-only study to identify genetic modifiers of breast cancer risk for BRCA1/BRCA2 mutation carriers; Table 5; PMID: 33597508)"
$ws.Range("C52").Value = "breast cancer (A case-only study to identify genetic modifiers of breast cancer risk for BRCA1/BRCA2 mutation carriers; Table 5; PMID: 33597508)"
$ws.Range("C53").Value = "breast cancer (A case-only study to identify genetic modifiers of breast cancer risk for BRCA1/BRCA2 mutation carriers; Table 5; PMID: 33597508)"
$ws.Range("C54").Value = "breast cancer (A case-only study to identify genetic modifiers of breast cancer risk for BRCA1/BRCA2 mutation carriers; Table 5; PMID: 33597508)"
$ws.Range("C55").Value = "breast cancer (A case-only study to identify genetic modifiers of breast cancer risk for BRCA1/BRCA2 mutation carriers; Table 5; PMID: 33597508)"

# rows 56-72: notes column (Table 6) - filled before rows 56-60 rs_id column
$ws.Range("C56").Value = "breast cancer (A case-only study to identify genetic modifiers of breast cancer risk for BRCA1/BRCA2 mutation carriers; Table 6; PMID: 33597508)"
$ws.Range("C57").Value = "breast cancer (A case-only study to identify genetic modifiers of breast cancer risk for BRCA1/BRCA2 mutation carriers; Table 6; PMID: 33597508)"
$ws.Range("C58").Value = "breast cancer (A case-only study to identify genetic modifiers of breast cancer risk for BRCA1/BRCA2 mutation carriers; Table 6; PMID: 33597508)"
$ws.Range("C59").Value = "breast cancer (A case-only study to identify genetic modifiers of breast cancer risk for BRCA1/BRCA2 mutation carriers; Table 6; PMID: 33597508)"
$ws.Range("C60").Value = "breast cancer (A case-only study to identify genetic modifiers of breast cancer risk for BRCA1/BRCA2 mutation carriers; Table 6; PMID: 33597508)"
$ws.Range("C61").Value = "breast cancer (A case-only study to identify genetic modifiers of breast cancer risk for BRCA1/BRCA2 mutation carriers; Table 6; PMID: 33597508)"
$ws.Range("C62").Value = "breast cancer (A case-only study to identify genetic modifiers of breast cancer risk for BRCA1/BRCA2 mutation carriers; Table 6; PMID: 33597508)"
$ws.Range("C63").Value = "breast cancer (A case-only study to identify genetic modifiers of breast cancer risk for BRCA1/BRCA2 mutation carriers; Table 6; PMID: 33597508)"
$ws.Range("C64").Value = "breast cancer (A case-only study to identify genetic modifiers of breast cancer risk for BRCA1/BRCA2 mutation carriers; Table 6; PMID: 33597508)"
$ws.Range("C65").Value = "breast cancer (A case-only study to identify genetic modifiers of breast cancer risk for BRCA1/BRCA2 mutation carriers; Table 6; PMID: 33597508)"
$ws.Range("C66").Value = "breast cancer (A case-only study to identify genetic modifiers of breast cancer risk for BRCA1/BRCA2 mutation carriers; Table 6; PMID: 33597508)"
$ws.Range("C67").Value = "breast cancer (A case-only study to identify genetic modifiers of breast cancer risk for BRCA1/BRCA2 mutation carriers; Table 6; PMID: 33597508)"
$ws.Range("C68").Value = "breast cancer (A case-only study to identify genetic modifiers of breast cancer risk for BRCA1/BRCA2 mutation carriers; Table 6; PMID: 33597508)"
$ws.Range("C69").Value = "breast cancer (A case-only study to identify genetic modifiers of breast cancer risk for BRCA1/BRCA2 mutation carriers; Table 6; PMID: 33597508)"
$ws.Range("C70").Value = "breast cancer (A case-only study to identify genetic modifiers of breast cancer risk for BRCA1/BRCA2 mutation carriers; Table 6; PMID: 33597508)"
$ws.Range("C71").Value = "breast cancer (A case-only study to identify genetic modifiers of breast cancer risk for BRCA1/BRCA2 mutation carriers; Table 6; PMID: 33597508)"
$ws.Range("C72").Value = "breast cancer (A case-only study to identify genetic modifiers of breast cancer risk for BRCA1/BRCA2 mutation carriers; Table 6; PMID: 33597508)"

# rows 56-60: rs_id column (Table 6)
$ws.Range("A56").Value = "rs12470785"
$ws.Range("A57").Value = "rs79183898"
$ws.Range("A58").Value = "rs71434801"
$ws.Range("A59").Value = "rs77197167"
$ws.Range("A60").Value = "rs114300732"

# row 61: secondary variant notation in the gene column
$ws.Range("B61").Value = "13:32231513:CAA:C"
# explicitly clear fill formatting on B61 so the cell gets its own cell format (applyFill)
$ws.Range("B61").Interior.ColorIndex = -4142

# rows 62-72: rs_id column
$ws.Range("A62").Value = "rs1623189"
$ws.Range("A63").Value = "rs736596"
$ws.Range("A64").Value = "rs77889880"
$ws.Range("A65").Value = "rs67776313"
$ws.Range("A66").Value = "rs71196514"
$ws.Range("A67").Value = "rs2555605"
$ws.Range("A68").Value = "rs74796280"
$ws.Range("A69").Value = "rs4943263"
$ws.Range("A70").Value = "rs2202781"
$ws.Range("A71").Value = "rs55675572"
$ws.Range("A72").Value = "rs17755120"

# row 61: rs_id (filled in last)
$ws.Range("A61").Value = "rs1198832427"

# restore the selection/active cell that was left after the edit
$ws.Range("B33").Select()

